$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.033.94"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "3.428.34"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'579.18"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("D6").Value = "'144.82"
$ws.Range("E6").Value = "  +2.23%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D11").Value = "'0.387"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "4.012.96"
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "'28.48"
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").Value = "'0.125"
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("D15").Value = "3.426.99"
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").Value = "62.058.31"
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("D18").Value = "'6.20"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").Value = "'13.97"
$ws.Range("E19").Value = "  +2.38%  "
$ws.Range("E20").Value = "  +2.82%  "
$ws.Range("D21").Value = "'393.04"
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("D22").Value = "'74.66"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("D26").Value = "'0.189"
$ws.Range("E26").Value = "  +2.50%  "
$ws.Range("E27").Value = "  +3.92%  "
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("D31").Value = "'1.42"
$ws.Range("E31").Value = "  +3.28%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("E33").Value = "  +1.67%  "
$ws.Range("D34").Value = "'5.30"
$ws.Range("E34").Value = "  +6.58%  "
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("D36").Value = "'167.96"
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("D37").Value = "3.460.33"
$ws.Range("E38").Value = "  +1.52%  "
$ws.Range("D39").Value = "'28.70"
$ws.Range("E39").Value = "  +7.48%  "
$ws.Range("D40").Value = "'0.0755"
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("E41").Value = "  +1.34%  "
$ws.Range("E43").Value = "  +1.50%  "
$ws.Range("E44").Value = "  +4.43%  "
$ws.Range("D45").Value = "2.517.54"
$ws.Range("E45").Value = "  +2.50%  "
$ws.Range("E46").Value = "  +0.57%  "
$ws.Range("D47").Value = "'6.64"
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("D48").Value = "'1.00"
$ws.Range("D49").Value = "'0.0264"
$ws.Range("E49").Value = "  +0.78%  "
$ws.Range("D50").Value = "'2.12"
$ws.Range("E50").Value = "  -0.78%  "
$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").Value = "'0.206"
$ws.Range("E51").Value = "  -0.25%  "
